$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update account numbers in column E for rows 2 and 3
$ws.Range("E2").Value = 1785991583
$ws.Range("E3").Value = 1785991583

# Row 2: change policy/patent data from SSA009 set to RPR007 set, and date
$ws.Range("L2").Value = "'03/02/2021"
$ws.Range("S2").Value = "RPR007"
$ws.Range("T2").Value = "ABC12RPR007"
$ws.Range("U2").Value = "ZAZ123RPR007"

# Row 3: change policy/patent data from SSA010 set to RPR008 set, and date
$ws.Range("L3").Value = "'03/02/2021"
$ws.Range("S3").Value = "RPR008"
$ws.Range("T3").Value = "ABC12RPR008"
$ws.Range("U3").Value = "ZAZ123RPR008"

# Row 4 keeps SSA011 data (content unchanged, only shared-string index shifts internally)
$ws.Range("S4").Value = "SSA011"
$ws.Range("T4").Value = "ABC12SSA011"
$ws.Range("U4").Value = "ZAZ123SSA011"

# Rows 5 & 6 keep their existing content too (unchanged logically)
$ws.Range("A5").Value = "ssurgwsoadev4-oci.opc.oracleoutsourcing.com"
$ws.Range("D5").Value = "gw"
$ws.Range("L5").Value = "'20/09/2020"
$ws.Range("M5").Value = "1.VUELVE CLIENTE"
$ws.Range("S5").Value = "RGR009"
$ws.Range("T5").Value = "ABC12RGR009"
$ws.Range("U5").Value = "ZAZ123RGR009"

$ws.Range("A6").Value = "ssurgwsoadev4-oci.opc.oracleoutsourcing.com"
$ws.Range("D6").Value = "gw"
$ws.Range("L6").Value = "'20/09/2020"
$ws.Range("S6").Value = "RGR010"
$ws.Range("T6").Value = "ABC12RGR010"
$ws.Range("U6").Value = "ZAZ123RGR010"

# Update view: scroll to show column L as leftmost, select S3
$ws.Application.ActiveWindow.ScrollColumn = 12
$ws.Range("S3").Select()

$wb.Save()
